$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item(1).Name = "Actuarial Firm Summary16_21"
$wb.Worksheets.Item(2).Name = "Actuarial Firm Summary 3_2021"

$ws1 = $wb.Worksheets.Item(1)

# Sheet1 header (unchanged names, kept for completeness)
$ws1.Cells.Item(1,1).Value = "actuarial_firm_name"
$ws1.Cells.Item(1,2).Value = "fy"
$ws1.Cells.Item(1,3).Value = "AAL"
$ws1.Cells.Item(1,4).Value = "AAL_percent"
$ws1.Cells.Item(1,5).Value = "UAL"

# Sheet1 data rows (firm name, fy, AAL, AAL_percent, UAL)
$ws1.Cells.Item(2,1).Value = "Gabriel, Roeder, Smith & Company (GRS)"
$ws1.Cells.Item(2,2).Value = 2021
$ws1.Cells.Item(2,3).Value = 1448217745970
$ws1.Cells.Item(2,4).Value = 0.132657780249413
$ws1.Cells.Item(2,5).Value = 261234414860
$ws1.Cells.Item(3,1).Value = "Gabriel, Roeder, Smith & Company (GRS)"
$ws1.Cells.Item(3,2).Value = 2016
$ws1.Cells.Item(3,3).Value = 1145910773040
$ws1.Cells.Item(3,4).Value = 0.104966245537585
$ws1.Cells.Item(3,5).Value = 353013780172
$ws1.Cells.Item(4,1).Value = "Cavanaugh Macdonald Consulting"
$ws1.Cells.Item(4,2).Value = 2021
$ws1.Cells.Item(4,3).Value = 901490366670
$ws1.Cells.Item(4,4).Value = 0.0825771616812854
$ws1.Cells.Item(4,5).Value = 154718673140
$ws1.Cells.Item(5,1).Value = "Milliman"
$ws1.Cells.Item(5,2).Value = 2021
$ws1.Cells.Item(5,3).Value = 791425215550
$ws1.Cells.Item(5,4).Value = 0.0724951151996524
$ws1.Cells.Item(5,5).Value = 49509898210
$ws1.Cells.Item(6,1).Value = "Milliman"
$ws1.Cells.Item(6,2).Value = 2016
$ws1.Cells.Item(6,3).Value = 684066156000
$ws1.Cells.Item(6,4).Value = 0.062660948639272
$ws1.Cells.Item(6,5).Value = 177510694110
$ws1.Cells.Item(7,1).Value = "Cavanaugh Macdonald Consulting"
$ws1.Cells.Item(7,2).Value = 2016
$ws1.Cells.Item(7,3).Value = 668830779530
$ws1.Cells.Item(7,4).Value = 0.0612653772692909
$ws1.Cells.Item(7,5).Value = 235119454970
$ws1.Cells.Item(8,1).Value = "Segal"
$ws1.Cells.Item(8,2).Value = 2021
$ws1.Cells.Item(8,3).Value = 649806185770
$ws1.Cells.Item(8,4).Value = 0.0595227108882364
$ws1.Cells.Item(8,5).Value = 160198525120
$ws1.Cells.Item(9,1).Value = "CalPERS"
$ws1.Cells.Item(9,2).Value = 2021
$ws1.Cells.Item(9,3).Value = 587976000000
$ws1.Cells.Item(9,4).Value = 0.0538590216954464
$ws1.Cells.Item(9,5).Value = 110653248000
$ws1.Cells.Item(10,1).Value = "Segal"
$ws1.Cells.Item(10,2).Value = 2016
$ws1.Cells.Item(10,3).Value = 510644048790
$ws1.Cells.Item(10,4).Value = 0.0467753597126944
$ws1.Cells.Item(10,5).Value = 193838029190
$ws1.Cells.Item(11,1).Value = "Scott Terando"
$ws1.Cells.Item(11,2).Value = 2016
$ws1.Cells.Item(11,3).Value = 436703008000
$ws1.Cells.Item(11,4).Value = 0.0400023075471426
$ws1.Cells.Item(11,5).Value = 137999008000
$ws1.Cells.Item(12,1).Value = "Conduent (formerly Buck and/or Xerox)"
$ws1.Cells.Item(12,2).Value = 2016
$ws1.Cells.Item(12,3).Value = 385162844740
$ws1.Cells.Item(12,4).Value = 0.035281191768255
$ws1.Cells.Item(12,5).Value = 132914565260
$ws1.Cells.Item(13,1).Value = "Cheiron"
$ws1.Cells.Item(13,2).Value = 2021
$ws1.Cells.Item(13,3).Value = 357816501213
$ws1.Cells.Item(13,4).Value = 0.0327762471628429
$ws1.Cells.Item(13,5).Value = 97309179763
$ws1.Cells.Item(14,1).Value = "New York City Office of the Actuary"
$ws1.Cells.Item(14,2).Value = 2021
$ws1.Cells.Item(14,3).Value = 248994048500
$ws1.Cells.Item(14,4).Value = 0.0228080327431707
$ws1.Cells.Item(14,5).Value = 7493808500
$ws1.Cells.Item(15,1).Value = "New York State and Local Retirement Systems' Actuary"
$ws1.Cells.Item(15,2).Value = 2021
$ws1.Cells.Item(15,3).Value = 231904000000
$ws1.Cells.Item(15,4).Value = 0.0212425720901207
$ws1.Cells.Item(15,5).Value = -28177076000
$ws1.Cells.Item(16,1).Value = "Internal Actuarial Services"
$ws1.Cells.Item(16,2).Value = 2016
$ws1.Cells.Item(16,3).Value = 221998704000
$ws1.Cells.Item(16,4).Value = 0.0203352398994125
$ws1.Cells.Item(16,5).Value = 75081150000
$ws1.Cells.Item(17,1).Value = "New York State And Local Retirement Systems' Actuary"
$ws1.Cells.Item(17,2).Value = 2016
$ws1.Cells.Item(17,3).Value = 202866000000
$ws1.Cells.Item(17,4).Value = 0.0185826705431317
$ws1.Cells.Item(17,5).Value = 19225796000
$ws1.Cells.Item(18,1).Value = "Conduent (formerly Buck and/or Xerox)"
$ws1.Cells.Item(18,2).Value = 2021
$ws1.Cells.Item(18,3).Value = 161404206940
$ws1.Cells.Item(18,4).Value = 0.0147847406753299
$ws1.Cells.Item(18,5).Value = 45675353140
$ws1.Cells.Item(19,1).Value = "Nystrs Office Of The Actuary"
$ws1.Cells.Item(19,2).Value = 2021
$ws1.Cells.Item(19,3).Value = 131077400000
$ws1.Cells.Item(19,4).Value = 0.0120067834918138
$ws1.Cells.Item(19,5).Value = -17071064000
$ws1.Cells.Item(20,1).Value = "Nystrs Office Of The Actuary"
$ws1.Cells.Item(20,2).Value = 2016
$ws1.Cells.Item(20,3).Value = 109305104000
$ws1.Cells.Item(20,4).Value = 0.0100124256224047
$ws1.Cells.Item(20,5).Value = 1798960000
$ws1.Cells.Item(21,1).Value = "Internal Actuarial Services"
$ws1.Cells.Item(21,2).Value = 2021
$ws1.Cells.Item(21,3).Value = 104534296000
$ws1.Cells.Item(21,4).Value = 0.00957541620097116
$ws1.Cells.Item(21,5).Value = 31406700000
$ws1.Cells.Item(22,1).Value = "Foster & Foster"
$ws1.Cells.Item(22,2).Value = 2021
$ws1.Cells.Item(22,3).Value = 98364489600
$ws1.Cells.Item(22,4).Value = 0.00901025752654515
$ws1.Cells.Item(22,5).Value = 28795731600
$ws1.Cells.Item(23,1).Value = "Office of The State Actuary - Washington"
$ws1.Cells.Item(23,2).Value = 2021
$ws1.Cells.Item(23,3).Value = 97128000000
$ws1.Cells.Item(23,4).Value = 0.00889699419574152
$ws1.Cells.Item(23,5).Value = -13019696000
$ws1.Cells.Item(24,1).Value = "Cheiron"
$ws1.Cells.Item(24,2).Value = 2016
$ws1.Cells.Item(24,3).Value = 94665352398
$ws1.Cells.Item(24,4).Value = 0.00867141391589276
$ws1.Cells.Item(24,5).Value = 25489705898
$ws1.Cells.Item(25,1).Value = "Perac"
$ws1.Cells.Item(25,2).Value = 2016
$ws1.Cells.Item(25,3).Value = 83529088000
$ws1.Cells.Item(25,4).Value = 0.00765132414042895
$ws1.Cells.Item(25,5).Value = 34615860000
$ws1.Cells.Item(26,1).Value = "Pwc"
$ws1.Cells.Item(26,2).Value = 2016
$ws1.Cells.Item(26,3).Value = 64963232000
$ws1.Cells.Item(26,4).Value = 0.00595067846594813
$ws1.Cells.Item(26,5).Value = 12781431200
$ws1.Cells.Item(27,1).Value = "Office of The State Actuary - Washington"
$ws1.Cells.Item(27,2).Value = 2016
$ws1.Cells.Item(27,3).Value = 61139000000
$ws1.Cells.Item(27,4).Value = 0.00560037608242156
$ws1.Cells.Item(27,5).Value = 5427847000
$ws1.Cells.Item(28,1).Value = "Pwc"
$ws1.Cells.Item(28,2).Value = 2021
$ws1.Cells.Item(28,3).Value = 59395150000
$ws1.Cells.Item(28,4).Value = 0.00544063817648049
$ws1.Cells.Item(28,5).Value = 2321468000
$ws1.Cells.Item(29,1).Value = "Korn Ferry Hay Group"
$ws1.Cells.Item(29,2).Value = 2021
$ws1.Cells.Item(29,3).Value = 52926848000
$ws1.Cells.Item(29,4).Value = 0.00484813709182619
$ws1.Cells.Item(29,5).Value = 12695832000
$ws1.Cells.Item(30,1).Value = "Korn Ferry Hay Group"
$ws1.Cells.Item(30,2).Value = 2016
$ws1.Cells.Item(30,3).Value = 47518964000
$ws1.Cells.Item(30,4).Value = 0.00435277105361637
$ws1.Cells.Item(30,5).Value = 21131076000
$ws1.Cells.Item(31,1).Value = "Bryan, Pendleton, Swats & Mcallister, Llc"
$ws1.Cells.Item(31,2).Value = 2016
$ws1.Cells.Item(31,3).Value = 45894705000
$ws1.Cells.Item(31,4).Value = 0.00420398776872034
$ws1.Cells.Item(31,5).Value = 2588496000
$ws1.Cells.Item(32,1).Value = "USI Consulting Group"
$ws1.Cells.Item(32,2).Value = 2021
$ws1.Cells.Item(32,3).Value = 45250140000
$ws1.Cells.Item(32,4).Value = 0.00414494515419335
$ws1.Cells.Item(32,5).Value = 13589248000
$ws1.Cells.Item(33,1).Value = "Lousiana Legislative Auditor"
$ws1.Cells.Item(33,2).Value = 2016
$ws1.Cells.Item(33,3).Value = 29272400000
$ws1.Cells.Item(33,4).Value = 0.0026813727544624
$ws1.Cells.Item(33,5).Value = 11734450000
$ws1.Cells.Item(34,1).Value = "Nyhart"
$ws1.Cells.Item(34,2).Value = 2016
$ws1.Cells.Item(34,3).Value = 24886090100
$ws1.Cells.Item(34,4).Value = 0.00227958363370398
$ws1.Cells.Item(34,5).Value = 12787232600
$ws1.Cells.Item(35,1).Value = "Foster & Foster"
$ws1.Cells.Item(35,2).Value = 2016
$ws1.Cells.Item(35,3).Value = 22339367550
$ws1.Cells.Item(35,4).Value = 0.00204630202854879
$ws1.Cells.Item(35,5).Value = 9039656310
$ws1.Cells.Item(36,1).Value = "Actuary South Dakota Retirement System"
$ws1.Cells.Item(36,2).Value = 2021
$ws1.Cells.Item(36,3).Value = 13865352000
$ws1.Cells.Item(36,4).Value = 0.00127007614967788
$ws1.Cells.Item(36,5).Value = -766847000
$ws1.Cells.Item(37,1).Value = "Cbiz"
$ws1.Cells.Item(37,2).Value = 2016
$ws1.Cells.Item(37,3).Value = 13679861000
$ws1.Cells.Item(37,4).Value = 0.00125308504154879
$ws1.Cells.Item(37,5).Value = 5206363000
$ws1.Cells.Item(38,1).Value = "G. S. Curran & Company"
$ws1.Cells.Item(38,2).Value = 2021
$ws1.Cells.Item(38,3).Value = 12346287600
$ws1.Cells.Item(38,4).Value = 0.001130928765301
$ws1.Cells.Item(38,5).Value = 665358900
$ws1.Cells.Item(39,1).Value = "South Dakota Retirement System"
$ws1.Cells.Item(39,2).Value = 2016
$ws1.Cells.Item(39,3).Value = 10851252000
$ws1.Cells.Item(39,4).Value = 0.00099398243617215
$ws1.Cells.Item(39,5).Value = 337790000
$ws1.Cells.Item(40,1).Value = "Bolton"
$ws1.Cells.Item(40,2).Value = 2021
$ws1.Cells.Item(40,3).Value = 10713525520
$ws1.Cells.Item(40,4).Value = 0.00098136659220172
$ws1.Cells.Item(40,5).Value = -693364263
$ws1.Cells.Item(41,1).Value = "G. S. Curran & Company"
$ws1.Cells.Item(41,2).Value = 2016
$ws1.Cells.Item(41,3).Value = 10297901900
$ws1.Cells.Item(41,4).Value = 0.00094329517165612
$ws1.Cells.Item(41,5).Value = 2302482390
$ws1.Cells.Item(42,1).Value = "Definiti"
$ws1.Cells.Item(42,2).Value = 2021
$ws1.Cells.Item(42,3).Value = 4179376800
$ws1.Cells.Item(42,4).Value = 0.000382833900949436
$ws1.Cells.Item(42,5).Value = 658566500
$ws1.Cells.Item(43,1).Value = "Findley"
$ws1.Cells.Item(43,2).Value = 2021
$ws1.Cells.Item(43,3).Value = 3793100000
$ws1.Cells.Item(43,4).Value = 0.000347450670083469
$ws1.Cells.Item(43,5).Value = -411732000
$ws1.Cells.Item(44,1).Value = "Aon"
$ws1.Cells.Item(44,2).Value = 2021
$ws1.Cells.Item(44,3).Value = 3635244300
$ws1.Cells.Item(44,4).Value = 0.00033299097517917
$ws1.Cells.Item(44,5).Value = 363093300
$ws1.Cells.Item(45,1).Value = "SilverStone Group"
$ws1.Cells.Item(45,2).Value = 2021
$ws1.Cells.Item(45,3).Value = 3518184500
$ws1.Cells.Item(45,4).Value = 0.000322268213862612
$ws1.Cells.Item(45,5).Value = 224574200
$ws1.Cells.Item(46,1).Value = "Southern Actuarial Services"
$ws1.Cells.Item(46,2).Value = 2021
$ws1.Cells.Item(46,3).Value = 2958058800
$ws1.Cells.Item(46,4).Value = 0.00027096029954557
$ws1.Cells.Item(46,5).Value = 527705800
$ws1.Cells.Item(47,1).Value = "Bps&M"
$ws1.Cells.Item(47,2).Value = 2016
$ws1.Cells.Item(47,3).Value = 2904694300
$ws1.Cells.Item(47,4).Value = 0.000266072073217852
$ws1.Cells.Item(47,5).Value = 216467500
$ws1.Cells.Item(48,1).Value = "SilverStone Group"
$ws1.Cells.Item(48,2).Value = 2016
$ws1.Cells.Item(48,3).Value = 2867807300
$ws1.Cells.Item(48,4).Value = 0.000262693197662931
$ws1.Cells.Item(48,5).Value = 625260800
$ws1.Cells.Item(49,1).Value = "Southern Actuarial Services"
$ws1.Cells.Item(49,2).Value = 2016
$ws1.Cells.Item(49,3).Value = 2152521300
$ws1.Cells.Item(49,4).Value = 0.000197172489007392
$ws1.Cells.Item(49,5).Value = 589469300
$ws1.Cells.Item(50,1).Value = "Nyhart"
$ws1.Cells.Item(50,2).Value = 2021
$ws1.Cells.Item(50,3).Value = 2100049300
$ws1.Cells.Item(50,4).Value = 0.000192366016317345
$ws1.Cells.Item(50,5).Value = 316280700
$ws1.Cells.Item(51,1).Value = "Hooker & Holcombe"
$ws1.Cells.Item(51,2).Value = 2021
$ws1.Cells.Item(51,3).Value = 1918985910
$ws1.Cells.Item(51,4).Value = 0.000175780480427681
$ws1.Cells.Item(51,5).Value = 470194990
$ws1.Cells.Item(52,1).Value = "Hooker & Holcombe"
$ws1.Cells.Item(52,2).Value = 2016
$ws1.Cells.Item(52,3).Value = 1425378000
$ws1.Cells.Item(52,4).Value = 0.000130565643200083
$ws1.Cells.Item(52,5).Value = 425714000
$ws1.Cells.Item(53,1).Value = "Britton Bender Pc"
$ws1.Cells.Item(53,2).Value = 2021
$ws1.Cells.Item(53,3).Value = 1313297900
$ws1.Cells.Item(53,4).Value = 0.000120299025961407
$ws1.Cells.Item(53,5).Value = 9753400.00000014
$ws1.Cells.Item(54,1).Value = "Mockenhaupt"
$ws1.Cells.Item(54,2).Value = 2016
$ws1.Cells.Item(54,3).Value = 1299809870
$ws1.Cells.Item(54,4).Value = 0.000119063512776517
$ws1.Cells.Item(54,5).Value = 604379950
$ws1.Cells.Item(55,1).Value = "Boomershine Consulting Group"
$ws1.Cells.Item(55,2).Value = 2021
$ws1.Cells.Item(55,3).Value = 703817630
$ws1.Cells.Item(55,4).Value = 0.0000644701977696499
$ws1.Cells.Item(55,5).Value = 151461190
$ws1.Cells.Item(56,1).Value = "Dean Actuaries"
$ws1.Cells.Item(56,2).Value = 2021
$ws1.Cells.Item(56,3).Value = 696707380
$ws1.Cells.Item(56,4).Value = 0.0000638188937895384
$ws1.Cells.Item(56,5).Value = -54478500
$ws1.Cells.Item(57,1).Value = "Conefry & Company, Llc"
$ws1.Cells.Item(57,2).Value = 2016
$ws1.Cells.Item(57,3).Value = 609079630
$ws1.Cells.Item(57,4).Value = 0.0000557921292815089
$ws1.Cells.Item(57,5).Value = 255469630
$ws1.Cells.Item(58,1).Value = "Boomershine Consulting Group"
$ws1.Cells.Item(58,2).Value = 2016
$ws1.Cells.Item(58,3).Value = 561973000
$ws1.Cells.Item(58,4).Value = 0.0000514771283168958
$ws1.Cells.Item(58,5).Value = 161951970
$ws1.Cells.Item(59,1).Value = "Dean Actuaries"
$ws1.Cells.Item(59,2).Value = 2016
$ws1.Cells.Item(59,3).Value = 545044060
$ws1.Cells.Item(59,4).Value = 0.0000499264253175541
$ws1.Cells.Item(59,5).Value = 56796680
$ws1.Cells.Item(60,1).Value = "Arthur J. Gallagher & Co."
$ws1.Cells.Item(60,2).Value = 2021
$ws1.Cells.Item(60,3).Value = 130922410
$ws1.Cells.Item(60,4).Value = 0.000011992586297077
$ws1.Cells.Item(60,5).Value = -1212889.99999998
$ws1.Cells.Item(61,1).Value = "Arthur J. Gallagher & Co."
$ws1.Cells.Item(61,2).Value = 2016
$ws1.Cells.Item(61,3).Value = 96540609
$ws1.Cells.Item(61,4).Value = 0.00000884318876046405
$ws1.Cells.Item(61,5).Value = 9793851
$ws1.Cells.Item(62,1).Value = "Usi Consulting Group"
$ws1.Cells.Item(62,2).Value = 2016
$ws1.Cells.Item(62,3).Value = 71594594
$ws1.Cells.Item(62,4).Value = 0.00000655811596310509
$ws1.Cells.Item(62,5).Value = 13546852
$ws1.Cells.Item(63,1).Value = "Principal Financial Group"
$ws1.Cells.Item(63,2).Value = 2021
$ws1.Cells.Item(63,3).Value = 65543141
$ws1.Cells.Item(63,4).Value = 0.00000600379854467989
$ws1.Cells.Item(63,5).Value = 1679539
$ws1.Cells.Item(64,1).Value = "McGriff Employee Benefit Solutions"
$ws1.Cells.Item(64,2).Value = 2021
$ws1.Cells.Item(64,3).Value = 63068992
$ws1.Cells.Item(64,4).Value = 0.00000577716472855684
$ws1.Cells.Item(64,5).Value = -5800328.00000001
$ws1.Cells.Item(65,1).Value = "Usi Consulting Group"
$ws1.Cells.Item(65,2).Value = 2021
$ws1.Cells.Item(65,3).Value = 61226578
$ws1.Cells.Item(65,4).Value = 0.00000560839828979404
$ws1.Cells.Item(65,5).Value = 6618937
$ws1.Cells.Item(66,1).Value = "Principal Financial Group"
$ws1.Cells.Item(66,2).Value = 2016
$ws1.Cells.Item(66,3).Value = 58870176
$ws1.Cells.Item(66,4).Value = 0.00000539255018299854
$ws1.Cells.Item(66,5).Value = 9890309
$ws1.Cells.Item(67,1).Value = "McGriff Employee Benefit Solutions"
$ws1.Cells.Item(67,2).Value = 2016
$ws1.Cells.Item(67,3).Value = 54127984
$ws1.Cells.Item(67,4).Value = 0.00000495816200761047
$ws1.Cells.Item(67,5).Value = 10294644

$ws2 = $wb.Worksheets.Item(2)

# Sheet2 header (AAL_percent sheet) plus new Year column
$ws2.Cells.Item(1,1).Value = "actuarial_firm_name"
$ws2.Cells.Item(1,2).Value = "AAL"
$ws2.Cells.Item(1,3).Value = "AAL_percent"
$ws2.Cells.Item(1,4).Value = "UAL"
$ws2.Cells.Item(1,5).Value = "Year"

# Sheet2 data rows (firm name, AAL, AAL_percent, UAL, Year)
$ws2.Cells.Item(2,1).Value = "CalPERS"
$ws2.Cells.Item(2,2).Value = 587976000000
$ws2.Cells.Item(2,3).Value = 0.0975121230195361
$ws2.Cells.Item(2,4).Value = 110653248000
$ws2.Cells.Item(2,5).Value = 2021
$ws2.Cells.Item(3,1).Value = "Cavanaugh Macdonald Consulting"
$ws2.Cells.Item(3,2).Value = 901490366670
$ws2.Cells.Item(3,3).Value = 0.149506509680075
$ws2.Cells.Item(3,4).Value = 154718673140
$ws2.Cells.Item(3,5).Value = 2021
$ws2.Cells.Item(4,1).Value = "Cheiron"
$ws2.Cells.Item(4,2).Value = 357816501213
$ws2.Cells.Item(4,3).Value = 0.0593416171488327
$ws2.Cells.Item(4,4).Value = 97309179763
$ws2.Cells.Item(4,5).Value = 2021
$ws2.Cells.Item(5,1).Value = "Conduent (formerly Buck and/or Xerox)"
$ws2.Cells.Item(5,2).Value = 161404206940
$ws2.Cells.Item(5,3).Value = 0.0267678729795161
$ws2.Cells.Item(5,4).Value = 45675353140
$ws2.Cells.Item(5,5).Value = 2021
$ws2.Cells.Item(6,1).Value = "Foster & Foster"
$ws2.Cells.Item(6,2).Value = 98364489600
$ws2.Cells.Item(6,3).Value = 0.0163131321867373
$ws2.Cells.Item(6,4).Value = 28795731600
$ws2.Cells.Item(6,5).Value = 2021
$ws2.Cells.Item(7,1).Value = "Gabriel, Roeder, Smith & Company (GRS)"
$ws2.Cells.Item(7,2).Value = 1448217745970
$ws2.Cells.Item(7,3).Value = 0.240177808284865
$ws2.Cells.Item(7,4).Value = 261234414860
$ws2.Cells.Item(7,5).Value = 2021
$ws2.Cells.Item(8,1).Value = "Internal Actuarial Services"
$ws2.Cells.Item(8,2).Value = 104534296000
$ws2.Cells.Item(8,3).Value = 0.0173363557888631
$ws2.Cells.Item(8,4).Value = 31406700000
$ws2.Cells.Item(8,5).Value = 2021
$ws2.Cells.Item(9,1).Value = "Korn Ferry Hay Group"
$ws2.Cells.Item(9,2).Value = 52926848000
$ws2.Cells.Item(9,3).Value = 0.00877758499192533
$ws2.Cells.Item(9,4).Value = 12695832000
$ws2.Cells.Item(9,5).Value = 2021
$ws2.Cells.Item(10,1).Value = "Milliman"
$ws2.Cells.Item(10,2).Value = 791425215550
$ws2.Cells.Item(10,3).Value = 0.131252896341814
$ws2.Cells.Item(10,4).Value = 49509898210
$ws2.Cells.Item(10,5).Value = 2021
$ws2.Cells.Item(11,1).Value = "New York City Office of the Actuary"
$ws2.Cells.Item(11,2).Value = 248994048500
$ws2.Cells.Item(11,3).Value = 0.04129409752858
$ws2.Cells.Item(11,4).Value = 7493808500
$ws2.Cells.Item(11,5).Value = 2021
$ws2.Cells.Item(12,1).Value = "New York State and Local Retirement Systems' Actuary"
$ws2.Cells.Item(12,2).Value = 231904000000
$ws2.Cells.Item(12,3).Value = 0.0384598204292735
$ws2.Cells.Item(12,4).Value = -28177076000
$ws2.Cells.Item(12,5).Value = 2021
$ws2.Cells.Item(13,1).Value = "Nystrs Office Of The Actuary"
$ws2.Cells.Item(13,2).Value = 131077400000
$ws2.Cells.Item(13,3).Value = 0.021738362711881
$ws2.Cells.Item(13,4).Value = -17071064000
$ws2.Cells.Item(13,5).Value = 2021
$ws2.Cells.Item(14,1).Value = "Office of The State Actuary - Washington"
$ws2.Cells.Item(14,2).Value = 97128000000
$ws2.Cells.Item(14,3).Value = 0.0161080681603356
$ws2.Cells.Item(14,4).Value = -13019696000
$ws2.Cells.Item(14,5).Value = 2021
$ws2.Cells.Item(15,1).Value = "Pwc"
$ws2.Cells.Item(15,2).Value = 59395150000
$ws2.Cells.Item(15,3).Value = 0.00985031221268181
$ws2.Cells.Item(15,4).Value = 2321468000
$ws2.Cells.Item(15,5).Value = 2021
$ws2.Cells.Item(16,1).Value = "Segal"
$ws2.Cells.Item(16,2).Value = 649806185770
$ws2.Cells.Item(16,3).Value = 0.107766270605705
$ws2.Cells.Item(16,4).Value = 160198525120
$ws2.Cells.Item(16,5).Value = 2021
$ws2.Cells.Item(17,1).Value = "Others"
$ws2.Cells.Item(17,2).Value = 107312888761
$ws2.Cells.Item(17,3).Value = 0.0177971679293788
$ws2.Cells.Item(17,4).Value = 15051100475
$ws2.Cells.Item(17,5).Value = 2021

Write-Host "edit complete"
